# BIOMD0000000012_annotations.xlsx — annotation fix
#
# The "sbml_type" column (column B) for the species rows (X, Y, Z, PX, PY, PZ
# in rows 9-14) was mislabeled "parameter". Correct it to "species", matching
# the updated runBioSimulations annotations described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9:B14").Value = "species"

# Leave the selection where the author left it after editing.
$ws.Range("B15").Select()
